# Runsheet edit: add an "exclude" / "keep" coding column (I) next to the
# existing comments column, reflecting which sessions are usable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("I1").Value = "exclude"

# Per-subject keep/exclude coding
$values = @(
    "keep",     # row 2  - SPEED_ACC_CHILD_GAZE_1
    "keep",     # row 3  - SPEED_ACC_CHILD_GAZE_2
    "keep",     # row 4  - SPEED_ACC_CHILD_GAZE_3
    "keep",     # row 5  - SPEED_ACC_CHILD_GAZE_4
    "keep",     # row 6  - SPEED_ACC_CHILD_GAZE_5
    "keep",     # row 7  - SPEED_ACC_CHILD_GAZE_6
    "exclude",  # row 8  - SPEED_ACC_CHILD_GAZE_7
    "keep",     # row 9  - SPEED_ACC_CHILD_GAZE_8
    "keep",     # row 10 - SPEED_ACC_CHILD_GAZE_9
    "keep",     # row 11 - SPEED_ACC_CHILD_GAZE_10
    "exclude",  # row 12 - SPEED_ACC_CHILD_GAZE_11
    "keep",     # row 13 - SPEED_ACC_CHILD_GAZE_12
    "keep",     # row 14 - SPEED_ACC_CHILD_GAZE_13
    "keep",     # row 15 - SPEED_ACC_CHILD_GAZE_14
    "keep",     # row 16 - SPEED_ACC_CHILD_GAZE_15
    "keep",     # row 17 - SPEED_ACC_CHILD_GAZE_16
    "keep",     # row 18 - SPEED_ACC_CHILD_GAZE_17
    "keep",     # row 19 - SPEED_ACC_CHILD_GAZE_18
    "keep",     # row 20 - SPEED_ACC_CHILD_GAZE_19
    "exclude",  # row 21 - SPEED_ACC_CHILD_GAZE_20
    "keep",     # row 22 - SPEED_ACC_CHILD_GAZE_21
    "exclude",  # row 23 - SPEED_ACC_CHILD_GAZE_22
    "keep",     # row 24 - SPEED_ACC_CHILD_GAZE_23
    "keep",     # row 25 - SPEED_ACC_CHILD_GAZE_24
    "keep",     # row 26 - SPEED_ACC_CHILD_GAZE_25
    "keep",     # row 27 - SPEED_ACC_CHILD_GAZE_26
    "keep",     # row 28 - SPEED_ACC_CHILD_GAZE_27
    "keep",     # row 29 - SPEED_ACC_CHILD_GAZE_28
    "keep",     # row 30 - SPEED_ACC_CHILD_GAZE_29
    "exclude",  # row 31 - SPEED_ACC_CHILD_GAZE_30
    "keep",     # row 32 - SPEED_ACC_CHILD_GAZE_31
    "keep",     # row 33 - SPEED_ACC_CHILD_GAZE_32
    "keep",     # row 34 - SPEED_ACC_CHILD_GAZE_33
    "keep",     # row 35 - SPEED_ACC_CHILD_GAZE_34
    "keep"      # row 36 - SPEED_ACC_CHILD_GAZE_35
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
}

# Scroll/selection state matching the author's final view
$ws.Range("A1").Select() | Out-Null
$ws.Range("I32:I36").Select() | Out-Null
